$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 533 }

$rng = $ws.Range("C2:C$lastRow")
$rng.Value = $newDate
